$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated x(Determined)/y(Determined)/Distance/Detected values per row (rows 4-43)
$changes = @{
    4  = @{ E = 99;  F = 264; G = 11.40175425099138 }
    5  = @{ E = 436; F = 278; G = 24.18677324489565 }
    6  = @{          F = 276; G = 34.0147027033899;  H = "Y" }
    7  = @{ E = 365; F = 280; G = 7.071067811865476 }
    8  = @{ E = 103; F = 260; G = 18 }
    9  = @{ E = 447; F = 274; G = 18.43908891458577 }
    10 = @{ E = 464; F = 285; G = 31.14482300479487 }
    11 = @{ E = 462; F = 283; G = 33.37663853655727 }
    12 = @{ E = 108; F = 266; G = 13.60147050873544 }
    13 = @{ E = 450; F = 284; G = 18.97366596101028 }
    14 = @{ E = 83;  F = 278; G = 11.31370849898476 }
    15 = @{ E = 107;          G = 22.3606797749979 }
    16 = @{ E = 476; F = 276; G = 9.055385138137417 }
    17 = @{ E = 461; F = 281; G = 22.02271554554524 }
    18 = @{ E = 181; F = 286; G = 15.23154621172782 }
    19 = @{ E = 461; F = 257; G = 19.23538406167134 }
    20 = @{ E = 447; F = 261; G = 21.37755832643195 }
    21 = @{ E = 450; F = 263; G = 25.17935662402834 }
    22 = @{ E = 473; F = 288; G = 14.42220510185596 }
    23 = @{ E = 474; F = 284; G = 7.211102550927978 }
    24 = @{ E = 89;  F = 256; G = 12.52996408614167 }
    25 = @{ E = 455; F = 281; G = 16.15549442140351 }
    26 = @{ E = 56;  F = 218; G = 379.0646382874562 }
    27 = @{ E = 459; F = 295; G = 16.76305461424021 }
    28 = @{ E = 445; F = 277; G = 20.09975124224178 }
    29 = @{ E = 103; F = 251; G = 15.03329637837291 }
    30 = @{ E = 476; F = 289; G = 12.72792206135786 }
    31 = @{ E = 475; F = 271; G = 15.03329637837291 }
    32 = @{ E = 483; F = 282; G = 16.97056274847714 }
    33 = @{ E = 470; F = 284; G = 20.8806130178211 }
    34 = @{ E = 415; F = 272; G = 33.30165161069343 }
    35 = @{ E = 477; F = 288; G = 13.15294643796591 }
    36 = @{ E = 451; F = 296; G = 21.02379604162864 }
    37 = @{ E = 265;          G = 126.2893503031827 }
    38 = @{ E = 95;  F = 274; G = 9.219544457292887;  H = "Y" }
    39 = @{ E = 452; F = 272; G = 19.69771560359221 }
    40 = @{ E = 476; F = 300; G = 17.20465053408525 }
    41 = @{ E = 485; F = 268; G = 8.602325267042627 }
    42 = @{ E = 118; F = 253; G = 25.94224354214569;  H = "Y" }
    43 = @{ E = 480; F = 284; G = 13.45362404707371 }
}

foreach ($row in $changes.Keys) {
    $vals = $changes[$row]
    if ($vals.ContainsKey("E")) { $ws.Range("E$row").Value = $vals["E"] }
    if ($vals.ContainsKey("F")) { $ws.Range("F$row").Value = $vals["F"] }
    if ($vals.ContainsKey("G")) { $ws.Range("G$row").Value = $vals["G"] }
    if ($vals.ContainsKey("H")) { $ws.Range("H$row").Value = $vals["H"] }
}
